$d = $word.ActiveDocument

# Replaces the text of a paragraph that exactly matches $oldText with $newText,
# while preserving the paragraph's existing run/formatting structure (including
# any empty leading <w:r/> runs) exactly as-is. Only the content of the <w:t>
# element in the paragraph's (single) text-bearing run is changed. This avoids
# the run-merging / smart-quote normalization that Find.Execute's built-in
# Replace performs.
function Replace-ParaText($oldText, $newText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        $trimmed = $t.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $oldText) {
            $full = $p.Range.WordOpenXML
            $body = $null
            if ($full -match '(?s)<w:body>(.*?)</w:body>') {
                $body = $matches[1]
            }
            $paraXml = $null
            if ($body -match '(?s)^(<w:p[^>]*>.*?</w:p>)') {
                $paraXml = $matches[1]
            }
            $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
            $newParaXml = $paraXml -replace '(?s)<w:t[^>]*>.*?</w:t>', ("<w:t xml:space=`"preserve`">" + $escaped + "</w:t>")
            $wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $p.Range.InsertXML($wrapped)
            return $true
        }
    }
    return $false
}

# Heading + bold footer line (same text appears twice in the doc)
Replace-ParaText "Play Musketeer Slot for Free: Game Review 2021" "Play Musketeer Slot Free - Exciting Slot Game Inspired by Alexandre Dumas' Novel" | Out-Null
Replace-ParaText "Play Musketeer Slot for Free: Game Review 2021" "Play Musketeer Slot Free - Exciting Slot Game Inspired by Alexandre Dumas' Novel" | Out-Null

# "What we like" bullets
Replace-ParaText "Expanding Wilds and Sticky Wilds" "Inspired by Alexandre Dumas' famous novel" | Out-Null
Replace-ParaText "Free Spins feature with Random Multipliers" "Expanding Wilds and Free Spins feature" | Out-Null
Replace-ParaText "Story-based game inspired by Alexandre Dumas' novel" "Random Multiplier for increased winnings" | Out-Null
Replace-ParaText "Playable on ADM licensed online casinos" "Available to play in demo mode" | Out-Null

# "What we don't like" bullets
Replace-ParaText "Low-paying standard card symbols" "Limited number of paylines" | Out-Null
Replace-ParaText "No progressive jackpot" "May not appeal to players unfamiliar with the novel" | Out-Null

# Italic summary line
Replace-ParaText "Read our expert review of Musketeer Slot by iSoftBet and play for free. Discover the game's features and winning potential, and find ADM licensed casinos to play for real money." "Read our review of Musketeer Slot, inspired by Alexandre Dumas' novel. Play for free and win big with expanding wilds and free spins." | Out-Null
